$wb = $excel.ActiveWorkbook

# --- "Rate laws" worksheet: update rate-law formula text and selection ---
$wsRateLaws = $wb.Worksheets.Item("Rate laws")
$wsRateLaws.Range("C2").Value = "growthRate * M[c]"
$wsRateLaws.Range("A3").Select() | Out-Null

# --- "Parameters" worksheet: incorporate 1/3600 h/s factor into growthRate parameter ---
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Range("E3").Value = 0.0000083713
$wsParameters.Range("E3").NumberFormat = "0.00E+00"
$wsParameters.Range("G3").Value = "ln(2)/23 h * 1 h / 3600 s = ln(2)/(23*3600) 1/s = 8.3713e-06 1/s"

# Make "Parameters" the active / selected sheet (matches Excel's tabSelected + activeTab change)
$wsParameters.Activate() | Out-Null
$wsParameters.Range("A4").Select() | Out-Null
